$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Result")

# Row 2 updates
$ws.Range("A2").Value = 338.90100000000001
$ws.Range("B2").Value = 250.501
$ws.Range("D2").Value = 0.029
$ws.Range("J2").Value = 150.07400000000001
$ws.Range("K2").Value = 169.79
$ws.Range("L2").Value = 240.47200000000001
$ws.Range("M2").Value = 209.846
$ws.Range("N2").Value = 144.875
$ws.Range("O2").Value = 129.816

# Row 3 updates
$ws.Range("A3").Value = 338.58
$ws.Range("B3").Value = 242.41300000000001
$ws.Range("J3").Value = 149.66200000000001
$ws.Range("K3").Value = 169.458
$ws.Range("L3").Value = 239.96100000000001
$ws.Range("M3").Value = 209.24100000000001
$ws.Range("N3").Value = 144.52699999999999
$ws.Range("O3").Value = 129.548
